$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 25 ("Update the changes on version"): clear out most of the row.
# C25, D25, G25, H25 are fully cleared (contents + formatting -> cell
# disappears from the sheet XML). E25 and F25 keep their formatting but lose
# their value (empty cell with style only).
# ---------------------------------------------------------------------------
$ws.Range("C25").Clear()
$ws.Range("D25").Clear()
$ws.Range("E25").ClearContents()
$ws.Range("F25").ClearContents()
$ws.Range("G25").Clear()
$ws.Range("H25").Clear()

# ---------------------------------------------------------------------------
# Row 27 used to only contain "Github update" in column A; that note is now
# relocated down to row 30, and row 27 gets brand-new content.
# ---------------------------------------------------------------------------
$row27 = $ws.Range("A27:H27")
$row27.Value = 1
$row27.HorizontalAlignment = -4108
$ws.Range("A27").Value = "remove 1 second timer stop from Ignition logic"

# Row 28 (new)
$row28 = $ws.Range("A28:H28")
$row28.Value = 1
$row28.HorizontalAlignment = -4108
$ws.Range("A28").Value = "No tag in when immobilizer activated, remove timer7 at tag"
$ws.Range("C28").Value = "n/a"
$ws.Range("E28").Value = "n/a"
$ws.Range("F28").Value = "n/a"

# Row 29 (new)
$row29 = $ws.Range("A29:H29")
$row29.Value = 1
$row29.HorizontalAlignment = -4108
$ws.Range("A29").Value = "add timer6 along with Tag in"
$ws.Range("C29").Value = "n/a"
$ws.Range("E29").Value = "n/a"
$ws.Range("F29").Value = "n/a"

# Row 30: relocated "Github update" note, column A only.
$ws.Range("A30").Value = "Github update"
$ws.Range("A30").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Sheet view bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("C29").Select()
$excel.ActiveWindow.ScrollRow = 10
